# Ordenamiento por base - add "n" / "nlog(n;10)" data table + chart
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: new header row for the second table -------------------------
# A14 used to hold the shared string "10f(n/10)+n"; it becomes "n" and a new
# header "nlog(n;10)" is added in B14. Both reuse the same look as the first
# table's header row (A1:B1).
$ws.Range("A14").Value2 = "n"
$ws.Range("B14").Value2 = "nlog(n;10)"
$ws.Range("A1:B1").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)

# --- Rows 15-25: n / n*log10(n) data --------------------------------------
$nValues = @(5, 10, 15, 20, 25, 30, 35, 40, 45, 50, 55)
for ($i = 0; $i -lt $nValues.Length; $i++) {
    $row = 15 + $i
    $ws.Cells.Item($row, 1).Value2 = $nValues[$i]
}

# Both columns of the new block use the same look: the Arial / theme-colour
# font already used for column A in the first table (cell A2).
$ws.Range("A2").Copy()
$ws.Range("A15:A25").PasteSpecial(-4122)
$ws.Range("B15:B25").PasteSpecial(-4122)

# B15:B18 share one formula group (relative reference).
$ws.Range("B15:B18").Formula = "=(A15)*LOG(A15,10)"

# B19:B25 repeat the same whole-range formula in every cell (legacy-style
# "array fill"). The leading "@" forces implicit intersection on the bare
# range operand so each row picks its own n instead of spilling/erroring.
for ($row = 19; $row -le 25; $row++) {
    $ws.Cells.Item($row, 2).Formula = "=@(A15:A25)*LOG(A15:A25,10)"
}

Write-Host "done"
